$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "seufzen",
    "reizen",
    "boxen",
    "schaden",
    "leisten",
    "jubeln",
    "sprengen",
    "sperren",
    "retten",
    "treiben",
    "freuen",
    "kosten",
    "husten",
    "sitzen",
    "dienen",
    "mögen",
    "reisen",
    "tropfen",
    "erben",
    "drehen",
    "wagen",
    "bleiben",
    "zünden",
    "machen",
    "gründen",
    "klappen",
    "heulen",
    "führen",
    "fallen",
    "wirken",
    "segnen",
    "baden",
    "spenden",
    "formen",
    "fischen",
    "strahlen",
    "treffen",
    "arten",
    "holen",
    "lesen",
    "streichen",
    "siegen",
    "brauchen",
    "nennen",
    "werfen",
    "schwingen",
    "zielen",
    "flehen",
    "kennen",
    "rufen",
    "fällen",
    "malen",
    "merken",
    "küssen",
    "scheiden",
    "greifen",
    "spinnen",
    "lügen",
    "loben",
    "läuten",
    "planen",
    "spüren",
    "wüten",
    "mauern",
    "stellen",
    "bieten",
    "folgen",
    "scheinen",
    "heben",
    "handeln",
    "atmen",
    "filmen",
    "sorgen",
    "schulden",
    "zögern",
    "kehren",
    "saufen",
    "fragen",
    "schließen",
    "lockern",
    "morden",
    "nähen",
    "pflanzen",
    "ärgern",
    "ändern",
    "fließen",
    "doppeln",
    "zeigen",
    "sichern",
    "sinken",
    "wundern",
    "wecken",
    "werden",
    "feiern",
    "zeichnen",
    "platzen",
    "dringen",
    "backen",
    "stammen",
    "hören",
    "schrecken",
    "wählen",
    "schreiten",
    "klingen",
    "bilden",
    "scheitern",
    "ehren",
    "warnen",
    "schlucken",
    "kümmern",
    "hoffen",
    "stoßen",
    "ahnen",
    "sparen",
    "münzen",
    "biegen",
    "jagen",
    "tragen",
    "wachsen",
    "hauen",
    "fahren",
    "garen",
    "bluten",
    "warten",
    "grüßen",
    "pfeifen",
    "töten",
    "liegen",
    "messen",
    "hacken",
    "decken",
    "schwören",
    "klettern",
    "schmecken",
    "fangen",
    "leihen",
    "flüstern",
    "liefern",
    "trennen",
    "stürzen",
    "streifen",
    "geben",
    "träumen",
    "lassen",
    "quälen",
    "rasen",
    "kichern",
    "kaufen",
    "räumen",
    "graben",
    "enden",
    "knarren",
    "schenken",
    "bellen",
    "lieben",
    "äußern",
    "wenden",
    "trösten",
    "stecken",
    "weigern",
    "gelten",
    "drohen",
    "spielen",
    "lösen",
    "sterben",
    "beißen",
    "irren",
    "weichen",
    "schwächen",
    "suchen",
    "parken",
    "bauen",
    "achten",
    "bergen",
    "wehtun",
    "flüchten",
    "gleichen",
    "runden",
    "trauen",
    "herrschen",
    "tollen",
    "wechseln",
    "dauern",
    "altern",
    "helfen",
    "betteln",
    "albern",
    "heißen",
    "heilen",
    "eignen",
    "bitten",
    "stehlen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $words[$i]
}

Write-Host "Done updating words"